$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'workout knee pad'
    'girls knee pads volleyball'
    'youth indoor volleyball'
    'boys compressions'
    'knee sleeve for bursitis'
    'stretch marks men'
    'mens below the knee shorts'
    'sport tights men'
    'joint compression'
    'softball sweats'
    'youth sport leggings for boys'
    'padded compression knee sleeve'
    'youth volleyball shorts'
    'knee hockey'
    'girls hockey compression shorts'
    'football pads adult'
    'padded football shorts'
    'basketball shorts for men xxl'
    'womens softball compression shorts'
    'basketball leg sleeves youth boys'
    'football youth tights'
    'paintball protective gear pants'
    'catcher pads'
    'size chart for men'
    'paintball shorts padded'
    'men spandex leggings'
    'calf compression sleeve for youth'
    'mens 3/4 shorts'
    'tendinitis knee'
    'basketball hopes'
    'mens volleyball shorts'
    'cycling pants men padded'
    'knee guards volleyball'
    'silicone calf pads'
    'work knees pads'
    'youth padded knee sleeves'
    'knee protector running'
    'sliding shorts youth boys'
    'padded leg sleeve'
    'snowboarding pants for men'
    'basketball sweat pants for men'
    'black baseball pants youth boys'
    '1/2x28 thread protector'
    'six one knee pads'
    'outdoor knee pads'
    'shorts with leggings men'
    'youth baseball compression'
    'hip guards for fall protection'
    'mens compression yoga pants'
    'basketball under shorts'
    'capris spandex'
    'men leg compression pants'
    'youth compression padded shorts'
    'youth volleyball knee pads for girls'
    'baseball short pants'
    'youth sports leggings boys'
    'youth athletic tights boys'
    'compression pant for men'
    'working knee pad'
    'mens long shorts below knee'
    'black youth football pants'
    'baseball pants youth xxl'
    'girls compression leggings'
    'easy knee pads'
    'compression knee sleeve youth'
    'basketball gear for boys youth'
    'men leggings black'
    'boy leggings for sports youth'
    'leggings men compression'
    'hockey pads youth'
    'under knee pad sleeves'
    'weightlifting floor'
    'anti sweat pads'
    'mens knee shorts'
    'leggings for sports men'
    'girls knee pads volleyball youth'
    'compression knee sleeve with pad'
    'running pad'
    'baseball stretch bands'
    'softball long pants'
    'elastic waisted pants for men'
    'mountain bike knee pads for men'
    'yoga pant for men'
    'compression basketball knee sleeve'
    'compression knee sleeves with pads'
    'knee work pad'
    'running compression pants'
    'professional knee pads construction'
    'compression pants for girls'
    'womans softball sliding shorts'
    'compression calf sleeve youth'
    'performance compression knee sleeve'
    'youth wrestling kneepads'
    'knee pads for work women'
    'baseball pants short'
    'men capri shorts'
    'volleyball knee pads girls'
    'sit pad hiking'
    'tactical pants knee pads'
    'training pants men'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
